$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (dSF) is column 6. Update rows 2-29 with re-pulled values.
$fValues = @{
    2  = -1
    3  = -1
    4  = 6
    5  = 2
    6  = 4
    7  = 8
    8  = 1
    9  = -3
    10 = -5
    11 = 2
    12 = 2
    14 = 2
    15 = 1
    17 = 1
    19 = 1
    20 = 1
    21 = 7
    22 = 1
    23 = 0
    24 = 3
    25 = 2
    27 = -4
    28 = -1
    29 = 3
}

foreach ($row in $fValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $fValues[$row]
}
